# prueba ok remesa abc
# Applies the edits described by the commit:
#  1. Row 7 ("VERDAD" agrupación) - the Agrupación label is replaced by its code "C009".
#  2. Row 20 (ABC) - Sociedades list loses the "C242" entry (C112,C200,C242 -> C112,C200).
#  3. A brand-new row for "C242" is inserted right after the "C235" row (row 25),
#     pushing every row below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. A7: VERDAD -> C009
$ws.Range("A7").Value = "C009"

# 2. B20: drop the trailing ",C242"
$ws.Range("B20").Value = "C112,C200"

# 3. Insert a new row below row 25 ("C235") for the "C242" entity that was split
#    out of the ABC agrupación, then fill it in.
# -4121 = xlShiftDown
$newRow = $ws.Rows.Item(26)
$newRow.Insert(-4121)

$ws.Range("A26").Value = "C242"
$ws.Range("B26").Value = "C242"
$ws.Range("C26").Value = "pcaballero@abc.es;gadanero@vocento.com, msantiago@abc.es, mafrias@abc.es"
